$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.401.99'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '3.646.44'
$ws.Range('E3').Value = '  +7.76%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('D7').Value = '3.636.27'
$ws.Range('E7').Value = '  +7.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.621'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.21%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  +0.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.609'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.64'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').Value = '4.235.78'
$ws.Range('E14').Value = '  +7.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '680.80'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '9.01'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.28%  '
$ws.Range('D17').Value = '3.650.19'
$ws.Range('E17').Value = '  +7.73%  '
$ws.Range('D18').Value = '71.578.55'
$ws.Range('E18').Value = '  +2.52%  '
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.59%  '
$ws.Range('E22').Value = '  +2.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.17'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +15.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '103.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('E27').Value = '  +4.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.43'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '579.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.32'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.83%  '
$ws.Range('E35').Value = '  +1.99%  '
$ws.Range('E36').Value = '  +2.45%  '
$ws.Range('D37').Value = '3.741.37'
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '35.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('E41').Value = '  +3.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.45'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0462'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.84%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('E47').Value = '  +4.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.133'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.11%  '
$ws.Range('E49').Value = '  +4.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '133.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.43%  '
